# Rename the "*img" sheets to "img*" and make the last one ("imge") the
# active/selected tab, matching the commit "Change names from *img to img*".

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# Move the active tab from "xbday" (index 3) to "imge" (index 16).
$wb.Worksheets.Item("imge").Activate()
